# Update gh-pages to output generated at 456a3b4
# Applies the latest scrape numbers (want-to-go counts / ticket prices)
# to both the "展览" (Exhibition) sheet and the consolidated "全部类型"
# (All types) sheet, which mirrors the same rows.

$wb = $excel.ActiveWorkbook

$exhibition = $wb.Worksheets.Item("展览")
$allTypes   = $wb.Worksheets.Item("全部类型")

# -- 展览 (sheet "展览") updates --------------------------------------
$exhibition.Range("F5").Value  = 2214
$exhibition.Range("G6").Value  = "已售罄"
$exhibition.Range("F7").Value  = 313
$exhibition.Range("F8").Value  = 1064
$exhibition.Range("F9").Value  = 1023
$exhibition.Range("F16").Value = 7869
$exhibition.Range("F25").Value = 1134
$exhibition.Range("F28").Value = 10
$exhibition.Range("F30").Value = 1665
$exhibition.Range("F34").Value = 46
$exhibition.Range("F36").Value = 277

# -- 全部类型 (sheet "全部类型") updates, same rows offset by +2 -------
$allTypes.Range("F7").Value  = 2214
$allTypes.Range("G8").Value  = "已售罄"
$allTypes.Range("F9").Value  = 313
$allTypes.Range("F10").Value = 1064
$allTypes.Range("F12").Value = 1023
$allTypes.Range("F19").Value = 7869
$allTypes.Range("F29").Value = 1134
$allTypes.Range("F32").Value = 10
$allTypes.Range("F34").Value = 1665
$allTypes.Range("F38").Value = 46
$allTypes.Range("F40").Value = 277

Write-Output "Applied scrape update 456a3b4 to 展览 and 全部类型 sheets"
